# Update "paises" (countries) COVID tracker sheet with refreshed data.
# The source data feed was re-pulled, which both updates several countries'
# metrics and re-sorts a few rows whose totals changed rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 24 de Junio de 2020 a las 16:22"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 2425855
$ws.Range("C4").Value = 1687
$ws.Range("D4").Value = 1020455
$ws.Range("E4").Value = 1281880
$ws.Range("G4").Value = 47
$ws.Range("H4").Value = 123520

# India (row 7)
$ws.Range("B7").Value = 461828
$ws.Range("C7").Value = 5713
$ws.Range("D7").Value = 261682
$ws.Range("E7").Value = 185569
$ws.Range("G7").Value = 94
$ws.Range("H7").Value = 14577

# Alemania (row 14)
$ws.Range("B14").Value = 192855
$ws.Range("C14").Value = 77
$ws.Range("E14").Value = 7569

# Argentina (row 33)
$ws.Range("B33").Value = 47216
$ws.Range("C33").Value = 13
$ws.Range("D33").Value = 13816
$ws.Range("E33").Value = 32315
$ws.Range("G33").Value = 7
$ws.Range("H33").Value = 1085

# Noruega (row 71)
$ws.Range("B71").Value = 8777
$ws.Range("C71").Value = 5
$ws.Range("E71").Value = 391

# Uzbekistan (row 76, no rank change)
$ws.Range("D76").Value = 4675
$ws.Range("E76").Value = 2153

# Rows 77-78 swap rank: Consejo Danes para los Refugiados overtakes Senegal
$ws.Range("A77").Value = "Consejo Danes para los Refugiados"
$ws.Range("B77").Value = 6213
$ws.Range("C77").Value = 186
$ws.Range("D77").Value = 870
$ws.Range("E77").Value = 5201
$ws.Range("G77").Value = 7
$ws.Range("H77").Value = 142

$ws.Range("A78").Value = "Senegal"
$ws.Range("B78").Value = 6129
$ws.Range("C78").Value = 95
$ws.Range("D78").Value = 4072
$ws.Range("E78").Value = 1964
$ws.Range("G78").Value = 4
$ws.Range("H78").Value = 93

# Rows 82-86 re-rank: Kenia, El Salvador, Guinea, Etiopia, Gabon
$ws.Range("A82").Value = "Kenia"
$ws.Range("B82").Value = 5206
$ws.Range("C82").Value = 254
$ws.Range("D82").Value = 1782
$ws.Range("E82").Value = 3296
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 128

$ws.Range("A83").Value = "El Salvador"
$ws.Range("B83").Value = 5150
$ws.Range("C83").Value = 177
$ws.Range("D83").Value = 2924
$ws.Range("E83").Value = 2107
$ws.Range("G83").Value = 6
$ws.Range("H83").Value = 119

$ws.Range("A84").Value = "Guinea"
$ws.Range("B84").Value = 5040
$ws.Range("D84").Value = 3685
$ws.Range("E84").Value = 1327
$ws.Range("H84").Value = 28

$ws.Range("A85").Value = "Etiopia"
$ws.Range("B85").Value = 5034
$ws.Range("C85").Value = 186
$ws.Range("D85").Value = 1486
$ws.Range("E85").Value = 3470
$ws.Range("G85").Value = 3
$ws.Range("H85").Value = 78

$ws.Range("A86").Value = "Gabon"
$ws.Range("B86").Value = 4849
$ws.Range("D86").Value = 2107
$ws.Range("E86").Value = 2703
$ws.Range("H86").Value = 39

# Cuba (row 103)
$ws.Range("B103").Value = 2319
$ws.Range("C103").Value = 1
$ws.Range("D103").Value = 2130
$ws.Range("E103").Value = 104

# Sri Lanka (row 108)
$ws.Range("B108").Value = 2001
$ws.Range("C108").Value = 10
$ws.Range("E108").Value = 428

# Libano (row 115)
$ws.Range("B115").Value = 1644
$ws.Range("C115").Value = 22
$ws.Range("D115").Value = 1103
$ws.Range("E115").Value = 508
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 33

# Birmania (row 161)
$ws.Range("B161").Value = 293
$ws.Range("C161").Value = 1
$ws.Range("D161").Value = 208
$ws.Range("E161").Value = 79

# Tied-value rows that just swap display order (rows 202-203, 208-209)
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("A209").Value = "Islas Malvinas"
